# Apply weekly update: shift the per-row market data (date, quality, volume,
# min/max/weighted prices, unit, origin, $/kg, kg/unit) so that each row now
# shows the following week's reading, per the source diff.
#
# Mapping of new-row -> data taken from old-row (a single 7-cycle among rows 2-8):
#   row2 <- old row7
#   row3 <- old row8
#   row4 <- old row5
#   row5 <- old row2
#   row6 <- old row4
#   row7 <- old row3
#   row8 <- old row6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (before) values for the columns that change, keyed by row.
# NOTE: use .Value2 (not .Value) for reading - in this runtime, .Value does not
# reliably return the underlying cell data.
$rows = @(2, 3, 4, 5, 6, 7, 8)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# new row -> source (old) row
$mapping = @{
    2 = 7
    3 = 8
    4 = 5
    5 = 2
    6 = 4
    7 = 3
    8 = 6
}

foreach ($newRow in $rows) {
    $srcRow = $mapping[$newRow]
    $data = $snapshot[$srcRow]

    $ws.Cells.Item($newRow, 4).Value2 = $data.D
    $ws.Cells.Item($newRow, 12).Value2 = $data.L
    $ws.Cells.Item($newRow, 13).Value2 = $data.M
    $ws.Cells.Item($newRow, 14).Value2 = $data.N
    $ws.Cells.Item($newRow, 15).Value2 = $data.O
    $ws.Cells.Item($newRow, 16).Value2 = $data.P
    $ws.Cells.Item($newRow, 17).Value2 = $data.Q
    $ws.Cells.Item($newRow, 18).Value2 = $data.R
    $ws.Cells.Item($newRow, 19).Value2 = $data.S
    $ws.Cells.Item($newRow, 20).Value2 = $data.T
}
